$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap rows 11 and 12 (columns F:V; A:E -- Indice/pais/torneio/temporada/data -- are unchanged) ---
$row11 = $ws.Range("F11:V11").Value2
$row12 = $ws.Range("F12:V12").Value2
$ws.Range("F11:V11").Value2 = $row12
$ws.Range("F12:V12").Value2 = $row11

# --- Swap rows 37 and 38 (columns F:V) ---
$row37 = $ws.Range("F37:V37").Value2
$row38 = $ws.Range("F38:V38").Value2
$ws.Range("F37:V37").Value2 = $row38
$ws.Range("F38:V38").Value2 = $row37

# --- Append two new match rows (46, 47), matching formatting of the last existing row (45) ---
$ws.Range("A45:V45").Copy()
$ws.Range("A46:V47").PasteSpecial(-4122)

# Row 46: Al Nasr 0 - 2 Al Bataeh
$ws.Range("A46").Value2 = 45
$ws.Range("B46").Value2 = "united-arab-emirates"
$ws.Range("C46").Value2 = "uae-league"
$ws.Range("D46").Value2 = "2023-2024"
$ws.Range("E46").Value2 = 45233.57291666666
$ws.Range("F46").Value2 = "Al Nasr"
$ws.Range("G46").Value2 = 0
$ws.Range("H46").Value2 = "Al Bataeh"
$ws.Range("I46").Value2 = 2
$ws.Range("J46").Value2 = 1.55
$ws.Range("K46").Value2 = "27/10/2023 17:43"
$ws.Range("L46").Value2 = 1.63
$ws.Range("M46").Value2 = "03/11/2023 13:37"
$ws.Range("N46").Value2 = 4.36
$ws.Range("O46").Value2 = "27/10/2023 17:43"
$ws.Range("P46").Value2 = 4.09
$ws.Range("Q46").Value2 = "03/11/2023 13:37"
$ws.Range("R46").Value2 = 4.86
$ws.Range("S46").Value2 = "27/10/2023 17:43"
$ws.Range("T46").Value2 = 5.2
$ws.Range("U46").Value2 = "03/11/2023 13:36"
$ws.Range("V46").Value2 = "https://www.betexplorer.com/football/united-arab-emirates/uae-league/al-nasr-al-bataeh/nwb6NoBc/"

# Row 47: Emirates Club 2 - 4 Hatta
$ws.Range("A47").Value2 = 46
$ws.Range("B47").Value2 = "united-arab-emirates"
$ws.Range("C47").Value2 = "uae-league"
$ws.Range("D47").Value2 = "2023-2024"
$ws.Range("E47").Value2 = 45233.57291666666
$ws.Range("F47").Value2 = "Emirates Club"
$ws.Range("G47").Value2 = 2
$ws.Range("H47").Value2 = "Hatta"
$ws.Range("I47").Value2 = 4
$ws.Range("J47").Value2 = 1.86
$ws.Range("K47").Value2 = "02/11/2023 16:23"
$ws.Range("L47").Value2 = 1.7
$ws.Range("M47").Value2 = "03/11/2023 13:17"
$ws.Range("N47").Value2 = 3.87
$ws.Range("O47").Value2 = "02/11/2023 16:23"
$ws.Range("P47").Value2 = 4.38
$ws.Range("Q47").Value2 = "03/11/2023 13:17"
$ws.Range("R47").Value2 = 3.58
$ws.Range("S47").Value2 = "02/11/2023 16:23"
$ws.Range("T47").Value2 = 4.33
$ws.Range("U47").Value2 = "03/11/2023 13:17"
$ws.Range("V47").Value2 = "https://www.betexplorer.com/football/united-arab-emirates/uae-league/emirates-club-hatta/SlF3M5Q3/"
